# The supplementary-material paragraph that described the
# "leave_out_term_comparison_results.xlsx" workbook (and the stray empty
# References-styled paragraph that trails it at the very end of the body)
# is removed. After this edit, the paragraph beginning "Similarly,
# comparing correlations across all measures..." becomes the last
# paragraph of the document body, immediately followed by the final
# section properties.

$d = $word.ActiveDocument

$markerText = "We also provide a sample of the leave-one-term-out results"

$deleteStartPos = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith($markerText)) {
        $deleteStartPos = $p.Range.Start
        break
    }
}

if ($deleteStartPos -ge 0) {
    $deleteEndPos = $d.Content.End
    $deleteRange = $d.Range($deleteStartPos, $deleteEndPos)
    $deleteRange.Delete()
}
